# Refresh cryptos price/volume snapshot (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.052.52'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '2.935.03'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '377.39'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.28'
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('E7').Value = '  -1.71%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -2.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.58'
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.139'
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '3.399.89'
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('E14').Value = '  -4.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.39'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = '2.924.57'
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.979'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '50.970.13'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.17'
$ws.Range('E19').Value = '  -10.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  -4.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.52'
$ws.Range('E21').Value = '  -5.59%  '
$ws.Range('D22').Value = '0.0₃0951'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.35'
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '262.00'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.90'
$ws.Range('E25').Value = '  +3.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.20'
$ws.Range('E26').Value = '  +8.12%  '
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.115'
$ws.Range('E28').Value = '  +8.39%  '
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.58'
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.11'
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.60'
$ws.Range('E34').Value = '  -2.17%  '
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.04'
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.97'
$ws.Range('E38').Value = '  -4.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.58'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.55'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('E42').Value = '  -5.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.78'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.29'
$ws.Range('E44').Value = '  -5.47%  '
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.275'
$ws.Range('E46').Value = '  -5.49%  '
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.006.05'
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.22'
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0346'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.484'
$ws.Range('E51').Value = '  +12.61%  '
